$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item("LP1912")

$ws.Cells.Item(2,1).Value = "Última actualización: 07:52:32"
$ws.Cells.Item(3,1).Value = "Total filas: 67"
$ws.Cells.Item(38,1).Value = "07:52:32"
$ws.Cells.Item(38,2).Value = "07:54"
$ws.Cells.Item(38,3).Value = "14_ABASTO"
$ws.Cells.Item(38,4).Value = 2
$ws.Cells.Item(38,5).Value = "LP1912"
$ws.Cells.Item(41,1).Value = "07:52:32"
$ws.Cells.Item(41,2).Value = "08:00"
$ws.Cells.Item(41,3).Value = "17_ROMERO"
$ws.Cells.Item(41,4).Value = 8
$ws.Cells.Item(41,5).Value = "LP1912"
$ws.Cells.Item(42,1).Value = "06:38:54"
$ws.Cells.Item(42,2).Value = "08:00"
$ws.Cells.Item(42,3).Value = "16_SANTA ANA"
$ws.Cells.Item(42,4).Value = 82
$ws.Cells.Item(42,5).Value = "LP1912"
$ws.Cells.Item(43,1).Value = "07:52:32"
$ws.Cells.Item(43,2).Value = "08:01"
$ws.Cells.Item(43,3).Value = "16_SANTA ANA"
$ws.Cells.Item(43,4).Value = 9
$ws.Cells.Item(43,5).Value = "LP1912"
$ws.Cells.Item(44,1).Value = "07:52:32"
$ws.Cells.Item(44,2).Value = "08:06"
$ws.Cells.Item(44,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(44,4).Value = 14
$ws.Cells.Item(44,5).Value = "LP1912"
$ws.Cells.Item(45,1).Value = "07:52:32"
$ws.Cells.Item(45,2).Value = "08:11"
$ws.Cells.Item(45,3).Value = "10_OLMOS"
$ws.Cells.Item(45,4).Value = 19
$ws.Cells.Item(45,5).Value = "LP1912"
$ws.Cells.Item(47,1).Value = "07:52:32"
$ws.Cells.Item(47,2).Value = "08:13"
$ws.Cells.Item(47,3).Value = "15X38_ABASTO"
$ws.Cells.Item(47,4).Value = 21
$ws.Cells.Item(47,5).Value = "LP1912"
$ws.Cells.Item(49,1).Value = "07:52:32"
$ws.Cells.Item(49,2).Value = "08:29"
$ws.Cells.Item(49,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(49,4).Value = 37
$ws.Cells.Item(49,5).Value = "LP1912"
$ws.Cells.Item(50,1).Value = "07:52:32"
$ws.Cells.Item(50,2).Value = "08:29"
$ws.Cells.Item(50,3).Value = "15_ABASTO"
$ws.Cells.Item(50,4).Value = 37
$ws.Cells.Item(50,5).Value = "LP1912"
$ws.Cells.Item(51,1).Value = "07:52:32"
$ws.Cells.Item(51,2).Value = "08:41"
$ws.Cells.Item(51,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(51,4).Value = 49
$ws.Cells.Item(51,5).Value = "LP1912"
$ws.Cells.Item(52,1).Value = "07:52:32"
$ws.Cells.Item(52,2).Value = "08:41"
$ws.Cells.Item(52,3).Value = "10_OLMOS"
$ws.Cells.Item(52,4).Value = 49
$ws.Cells.Item(52,5).Value = "LP1912"
$ws.Cells.Item(54,1).Value = "07:52:32"
$ws.Cells.Item(54,2).Value = "08:44"
$ws.Cells.Item(54,3).Value = "215C_EL PATO"
$ws.Cells.Item(54,4).Value = 52
$ws.Cells.Item(54,5).Value = "LP1912"
$ws.Cells.Item(56,1).Value = "07:52:32"
$ws.Cells.Item(56,2).Value = "08:47"
$ws.Cells.Item(56,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(56,4).Value = 55
$ws.Cells.Item(56,5).Value = "LP1912"
$ws.Cells.Item(57,1).Value = "07:15:48"
$ws.Cells.Item(57,2).Value = "08:51"
$ws.Cells.Item(57,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(57,4).Value = 96
$ws.Cells.Item(57,5).Value = "LP1912"
$ws.Cells.Item(58,1).Value = "06:56:24"
$ws.Cells.Item(58,2).Value = "08:52"
$ws.Cells.Item(58,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(58,4).Value = 116
$ws.Cells.Item(58,5).Value = "LP1912"
$ws.Cells.Item(59,1).Value = "07:52:32"
$ws.Cells.Item(59,2).Value = "08:53"
$ws.Cells.Item(59,3).Value = "215B_EL PATO"
$ws.Cells.Item(59,4).Value = 61
$ws.Cells.Item(59,5).Value = "LP1912"
$ws.Cells.Item(60,1).Value = "07:40:11"
$ws.Cells.Item(60,2).Value = "08:54"
$ws.Cells.Item(60,3).Value = "215B_EL PATO"
$ws.Cells.Item(60,4).Value = 74
$ws.Cells.Item(60,5).Value = "LP1912"
$ws.Cells.Item(61,1).Value = "07:15:48"
$ws.Cells.Item(61,2).Value = "08:57"
$ws.Cells.Item(61,3).Value = "215A_EL PATO"
$ws.Cells.Item(61,4).Value = 102
$ws.Cells.Item(61,5).Value = "LP1912"
$ws.Cells.Item(62,1).Value = "07:52:32"
$ws.Cells.Item(62,2).Value = "08:58"
$ws.Cells.Item(62,3).Value = "215A_EL PATO"
$ws.Cells.Item(62,4).Value = 66
$ws.Cells.Item(62,5).Value = "LP1912"
$ws.Cells.Item(63,1).Value = "07:52:32"
$ws.Cells.Item(63,2).Value = "09:06"
$ws.Cells.Item(63,3).Value = "16_SANTA ANA"
$ws.Cells.Item(63,4).Value = 74
$ws.Cells.Item(63,5).Value = "LP1912"
$ws.Cells.Item(64,1).Value = "07:52:32"
$ws.Cells.Item(64,2).Value = "09:12"
$ws.Cells.Item(64,3).Value = "27_EL RETIRO"
$ws.Cells.Item(64,4).Value = 80
$ws.Cells.Item(64,5).Value = "LP1912"
$ws.Cells.Item(65,1).Value = "07:40:11"
$ws.Cells.Item(65,2).Value = "09:14"
$ws.Cells.Item(65,3).Value = "27_EL RETIRO"
$ws.Cells.Item(65,4).Value = 94
$ws.Cells.Item(65,5).Value = "LP1912"
$ws.Cells.Item(66,1).Value = "07:52:32"
$ws.Cells.Item(66,2).Value = "09:17"
$ws.Cells.Item(66,3).Value = "14_ABASTO"
$ws.Cells.Item(66,4).Value = 85
$ws.Cells.Item(66,5).Value = "LP1912"
$ws.Cells.Item(67,1).Value = "07:40:11"
$ws.Cells.Item(67,2).Value = "09:18"
$ws.Cells.Item(67,3).Value = "14_ABASTO"
$ws.Cells.Item(67,4).Value = 98
$ws.Cells.Item(67,5).Value = "LP1912"
$ws.Cells.Item(68,1).Value = "07:52:32"
$ws.Cells.Item(68,2).Value = "09:18"
$ws.Cells.Item(68,3).Value = "15X38_ABASTO"
$ws.Cells.Item(68,4).Value = 86
$ws.Cells.Item(68,5).Value = "LP1912"
$ws.Cells.Item(69,1).Value = "07:40:11"
$ws.Cells.Item(69,2).Value = "09:31"
$ws.Cells.Item(69,3).Value = "16_SANTA ANA"
$ws.Cells.Item(69,4).Value = 111
$ws.Cells.Item(69,5).Value = "LP1912"
$ws.Cells.Item(70,1).Value = "07:52:32"
$ws.Cells.Item(70,2).Value = "09:39"
$ws.Cells.Item(70,3).Value = "15_ABASTO"
$ws.Cells.Item(70,4).Value = 107
$ws.Cells.Item(70,5).Value = "LP1912"
$ws.Cells.Item(71,1).Value = "07:52:32"
$ws.Cells.Item(71,2).Value = "09:41"
$ws.Cells.Item(71,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(71,4).Value = 109
$ws.Cells.Item(71,5).Value = "LP1912"
$ws.Cells.Item(72,1).Value = "07:52:32"
$ws.Cells.Item(72,2).Value = "09:43"
$ws.Cells.Item(72,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(72,4).Value = 111
$ws.Cells.Item(72,5).Value = "LP1912"

# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Cells.Item(2,1).Value = "Última actualización: 07:52:32"
$ws.Cells.Item(14,1).Value = "07:52:32"
$ws.Cells.Item(14,2).Value = "08:44"
$ws.Cells.Item(14,3).Value = "215C_EL PATO"
$ws.Cells.Item(14,4).Value = 52
$ws.Cells.Item(14,5).Value = "LP1912"
$ws.Cells.Item(15,1).Value = "07:52:32"
$ws.Cells.Item(15,2).Value = "08:53"
$ws.Cells.Item(15,3).Value = "215B_EL PATO"
$ws.Cells.Item(15,4).Value = 61
$ws.Cells.Item(15,5).Value = "LP1912"
$ws.Cells.Item(18,1).Value = "07:52:32"
$ws.Cells.Item(18,2).Value = "08:58"
$ws.Cells.Item(18,3).Value = "215A_EL PATO"
$ws.Cells.Item(18,4).Value = 66
$ws.Cells.Item(18,5).Value = "LP1912"

# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Cells.Item(2,1).Value = "Última actualización: 07:52:32"
$ws.Cells.Item(3,1).Value = "Total filas: 7"
$ws.Cells.Item(9,1).Value = "07:52:32"
$ws.Cells.Item(9,2).Value = "08:36"
$ws.Cells.Item(9,3).Value = "215A_LA PLATA"
$ws.Cells.Item(9,4).Value = 44
$ws.Cells.Item(9,5).Value = "L6173"
$ws.Cells.Item(11,1).Value = "07:52:32"
$ws.Cells.Item(11,2).Value = "08:51"
$ws.Cells.Item(11,3).Value = "215C_LA PLATA"
$ws.Cells.Item(11,4).Value = 59
$ws.Cells.Item(11,5).Value = "L6203"
$ws.Cells.Item(12,1).Value = "07:52:32"
$ws.Cells.Item(12,2).Value = "09:21"
$ws.Cells.Item(12,3).Value = "215A_LA PLATA"
$ws.Cells.Item(12,4).Value = 89
$ws.Cells.Item(12,5).Value = "L6173"
